$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 457 (Femacal de La Calera - Berenjena),
# which pushes all the subsequent rows (old 457..535) down by one (new 458..536).
$ws.Rows.Item(457).Insert()

# Populate the newly inserted row 457 with the new record's data.
$ws.Range("A457").Value = 3
$ws.Range("B457").Value = "Femacal de La Calera"
$ws.Range("C457").Value = "Coquimbo"
$ws.Range("D457").Value = 45218
$ws.Range("E457").Value = 5
$ws.Range("F457").Value = 100112001
$ws.Range("G457").Value = "Berenjena"
$ws.Range("H457").Value = "Sin especificar"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 45
$ws.Range("K457").Value = 8000
$ws.Range("L457").Value = 8000
$ws.Range("M457").Value = 8000
$ws.Range("N457").Value = "$/caja 60 unidades"
$ws.Range("O457").Value = "Región de Arica y Parinacota"
$ws.Range("P457").Value = 133
$ws.Range("Q457").Value = 60
$ws.Range("R457").Value = "Hortaliza"
